$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures -------------------------------------------------
# "VALOR MORA" total grows with the new worker's debt
$ws.Range("E11").Value = 190044
# "Cant. Trabajadores" (worker count) goes from 4 to 5
$ws.Range("C13").Value = 5

# --- Add a new worker row (row 20), pushing the old "last row" style down ---
# First, clone row 19's current ("last row") formatting onto the new row 20
# so the new last data row keeps the bordered/highlighted look.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Then, clone row 18's ("normal row") formatting onto row 19, since row 19 is
# no longer the last row in the table.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the new worker (row 20) with its data
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047421075"
$ws.Range("D20").Value = "JOSE CARLOS VERGARA PEREZ"
$ws.Range("E20").Value = "2201"
$ws.Range("F20").Value = 44680
$ws.Range("G20").Value = 877803

# --- Rework the signature block (rows 24-25 -> 25-26) -----------------------
# Insert a new row above the current row 25 ("NOMBRE DEL REPRESENTANTE LEGAL"
# / "FIRMA DEL REPRESENTANTE LEGAL"). This shifts that row's values, format
# and merged cells down to row 26 for free, and leaves a styled-but-blank
# duplicate behind at row 25.
$ws.Rows("25").Insert(-4121)

# Re-merge the new row 25 (the insert doesn't carry merges to the newly
# created row) and turn it into the signature underline.
$ws.Range("B25:C25").Merge()
$ws.Range("H25:J25").Merge()
$ws.Range("B25").Value = "___________________________________"
$ws.Range("H25").Value = "___________________________________"

# Remove the old "Observaciones" row (24) entirely - it no longer exists.
$ws.Range("B24:C24").UnMerge()
$ws.Range("H24:J24").UnMerge()
$ws.Range("B24:J24").Clear()
